$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 289, shifting existing rows 289-335 down to 290-336.
$ws.Rows(289).Insert()

# Populate the newly inserted row 289 with the new record's data.
$ws.Cells.Item(289, 1).Value = 6
$ws.Cells.Item(289, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(289, 3).Value = "Metropolitana"
$ws.Cells.Item(289, 4).Value = 45154
$ws.Cells.Item(289, 5).Value = 13
$ws.Cells.Item(289, 6).Value = 100112022
$ws.Cells.Item(289, 7).Value = "Arveja Verde"
$ws.Cells.Item(289, 8).Value = "Perfection"
$ws.Cells.Item(289, 9).Value = "Primera"
$ws.Cells.Item(289, 10).Value = 410
$ws.Cells.Item(289, 11).Value = 30000
$ws.Cells.Item(289, 12).Value = 32000
$ws.Cells.Item(289, 13).Value = 30878
$ws.Cells.Item(289, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(289, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(289, 16).Value = 1235
$ws.Cells.Item(289, 17).Value = 25
$ws.Cells.Item(289, 18).Value = "Hortaliza"
